$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, copying the style used by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save indicator values for rows 2-15
$saveValues = @(1, 0, 0, 1, 1, 0, 0, 1, 1, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
